$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cell H1 - copy formatting from the adjacent header cell (G1)
# then set its value/text.
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)
$ws.Range("H1").Value = "Save"

# New data cell H2 (plain number, no special style - like B2:G2)
$ws.Range("H2").Value = 1
